$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.221.27"
$ws.Range("E2").Value = "  -1.43%  "

$ws.Range("D3").Value = "2.929.13"
$ws.Range("E3").Value = "  -2.85%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "567.23"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.15%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "158.01"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.21%  "

$ws.Range("E7").Value = "  +0.11%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.514"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.60%  "

$ws.Range("D9").Value = "2.924.54"
$ws.Range("E9").Value = "  -2.86%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.70"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.68%  "

$ws.Range("E11").Value = "  -4.07%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.460"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.61%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000243"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.35%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.21"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.29%  "

$ws.Range("E15").Value = "  -0.84%  "

$ws.Range("D16").Value = "65.240.77"
$ws.Range("E16").Value = "  -1.30%  "

$ws.Range("D17").Value = "3.419.39"
$ws.Range("E17").Value = "  -2.68%  "

$ws.Range("E18").Value = "  +0.18%  "

$ws.Range("D19").Value = "2.930.63"
$ws.Range("E19").Value = "  -2.83%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.60"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +12.64%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "442.87"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.50%  "

$ws.Range("E22").Value = "  +0.78%  "

$ws.Range("E23").Value = "  -1.79%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "82.07"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.13%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.21"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.19%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.06"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.70%  "

$ws.Range("B27").Value = "Dai"
$ws.Range("C27").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.09%  "

$ws.Range("B28").Value = "RenderToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.01"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -6.26%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.03"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.47%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.35"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.65%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.57"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.76%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0000100"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.42%  "

$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.111"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.05%  "

$ws.Range("B34").Value = "EthereumClassic"
$ws.Range("C34").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "27.05"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.08%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.04%  "

$ws.Range("E36").Value = "  -2.31%  "

$ws.Range("E37").Value = "  -1.62%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "49.72"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.95%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "44.54"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.28%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.97"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -9.46%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.299"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.62%  "

$ws.Range("E42").Value = "  -2.58%  "

$ws.Range("E43").Value = "  -7.74%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.47"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.11%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "381.98"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.07%  "

$ws.Range("E46").Value = "  -1.08%  "

$ws.Range("D47").Value = "2.698.57"
$ws.Range("E47").Value = "  -3.56%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "133.75"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.37%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.19"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +4.60%  "

$ws.Range("B51").Value = "Stellar"
$ws.Range("C51").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.106"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.04%  "
